# dev of attack alliance village
# Adds a new "collectLevel" sheet (copied/derived from "vipLevel"), with an
# extra FLOAT_collectPercentPerHour column, tweaks the vipLevel expTo cap,
# drops the stale external workbook link, and restores the various sheet
# selections left behind by the edit session.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Break the now-unused external reference to PlayerVillageExp.xlsx
#    (removes xl/externalLinks/*, the <externalReferences> block, and the
#    matching relationship / content-type entries).
# ---------------------------------------------------------------------
foreach ($src in $wb.LinkSources(1)) {
    $wb.BreakLink($src, 1)
}

# ---------------------------------------------------------------------
# 2. resources sheet: just a leftover cursor move, no data changes.
# ---------------------------------------------------------------------
$wsResources = $wb.Worksheets.Item("resources")
$wsResources.Activate()
$wsResources.Range("H11").Select()

# ---------------------------------------------------------------------
# 3. vipLevel sheet: the expTo cap for the last row shrinks from
#    100,000,000 down to 10,000.
# ---------------------------------------------------------------------
$wsVip = $wb.Worksheets.Item("vipLevel")
$wsVip.Activate()
$wsVip.Range("C11").Value = 10000
$wsVip.Range("C12").Select()

# ---------------------------------------------------------------------
# 4. New collectLevel sheet: duplicate vipLevel (keeps styles/number
#    formats intact) directly after it, rename, then add the new
#    FLOAT_collectPercentPerHour column D.
# ---------------------------------------------------------------------
$wsVip.Copy($null, $wsVip)
$wsCollect = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsCollect.Name = "collectLevel"

$wsCollect.Range("D1").Value = "FLOAT_collectPercentPerHour"
$wsCollect.Range("D2").Value = 0.1
$wsCollect.Range("D3").Value = 0.15
$wsCollect.Range("D4").Value = 0.2
$wsCollect.Range("D5").Value = 0.25
$wsCollect.Range("D6").Value = 0.3
$wsCollect.Range("D7").Value = 0.35
$wsCollect.Range("D8").Value = 0.4
$wsCollect.Range("D9").Value = 0.45
$wsCollect.Range("D10").Value = 0.5
$wsCollect.Range("D11").Value = 0.55

$wsCollect.Activate()
$wsCollect.Range("D12").Select()

# ---------------------------------------------------------------------
# 5. Re-create the per-sheet Excel_BuiltIn__FilterDatabase_* defined
#    names so the new collectLevel sheet (localSheetId 7) gets its own
#    entries, ordered the way Excel emits them (new sheet first, then
#    the pre-existing sheets in their original order, then the
#    workbook-scoped entry).
# ---------------------------------------------------------------------
$i = $wb.Names.Count
while ($i -ge 1) {
    $n = $wb.Names.Item($i)
    if ($n.Name -like "*FilterDatabase*") {
        $n.Delete()
    }
    $i = $i - 1
}

$filterSheetOrder = @("collectLevel", "dragonMaterials", "houses", "materials", "playerLevel", "soldierMaterials", "vipLevel")
$filterBases = @("Excel_BuiltIn__FilterDatabase_2", "Excel_BuiltIn__FilterDatabase_6", "Excel_BuiltIn__FilterDatabase_7", "Excel_BuiltIn__FilterDatabase_8")
foreach ($base in $filterBases) {
    foreach ($sn in $filterSheetOrder) {
        $wsTarget = $wb.Worksheets.Item($sn)
        $wsTarget.Names.Add($base, "=#REF!")
    }
    $wb.Names.Add($base, "=#REF!")
}
